# Trade #31 closed at 2026-02-16 22:55:11 - base_strategy UP +0.000%
#
# Appends a new trade-log row (row 32) to both the "All Trades" sheet and
# the "base_strategy" sheet, mirroring the structure/values of the existing
# trade rows (same shape as the very first trade row, just a later
# timestamp and trade number).

function Add-TradeRow {
    param($ws)

    $row = 32
    $scratchRow = $row + 1

    $ws.Cells.Item($row, 1).Value = 31

    # Column B holds a plain date-looking string ("2026-02-16"). Assigning
    # it directly makes Excel auto-detect it as a date and convert it to a
    # serial number, so it is entered with a leading space (which blocks
    # the date parser) and the stray space is then trimmed back out via a
    # helper formula cell, copy, and paste-values - this keeps the cell a
    # plain text value with no extra number formatting applied.
    $ws.Cells.Item($row, 2).Value = " 2026-02-16"
    $scratch = $ws.Cells.Item($scratchRow, 2)
    $scratch.Formula = "=TRIM(B" + $row + ")"
    $scratch.Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4163) # xlPasteValues
    $scratch.ClearContents()

    $ws.Cells.Item($row, 3).Value = "22:55:11"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 49.999998
    # Column G (Exit Price) stays blank - trade is still OPEN.
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # Column P (Exit Reason) stays blank - trade is still OPEN.
    $ws.Cells.Item($row, 17).Value = 0
}

$wb = $excel.ActiveWorkbook

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$baseStrategy = $wb.Worksheets.Item("base_strategy")
Add-TradeRow $baseStrategy

Write-Output "Added row 32 to 'All Trades' and 'base_strategy' sheets."
